# Updated symbol list on Fri Jan 13 04:56:55 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the rows whose quotes changed, while keeping them stored as literal text
# (matching the source feed, which always renders these as text strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "289.72"
    "E2" = "1.01%"
    "D3" = "29.52"
    "E3" = "4.16%"
    "D4" = "5.129"
    "E4" = "4.31%"
    "E5" = "2.10%"
    "E6" = "1.55%"
    "E7" = "-1.06%"
    "D8" = "0.9163"
    "E8" = "0.46%"
    "D9" = "0.1588"
    "E9" = "1.22%"
    "D10" = "0.06679"
    "E10" = "-0.25%"
    "D11" = "0.07701"
    "E11" = "0.07%"
    "E12" = "-1.35%"
    "D13" = "0.08985"
    "E13" = "0.07%"
    "D14" = "0.001585"
    "E14" = "-1.29%"
    "D15" = "0.04524"
    "E15" = "1.35%"
    "D16" = "0.0006469"
    "D17" = "0.006261"
    "E17" = "3.91%"
    "D18" = "3.452"
    "E18" = "-0.53%"
    "D19" = "3.407"
    "D20" = "2.218"
    "E20" = "-1.03%"
    "E21" = "1.78%"
    "D22" = "0.1310"
    "E22" = "-2.95%"
    "D23" = "4.062"
    "D24" = "0.1566"
    "E24" = "3.02%"
    "D25" = "0.001188"
    "E25" = "0.24%"
    "D26" = "0.004130"
    "E26" = "-4.70%"
    "E27" = "5.91%"
    "E28" = "-1.14%"
    "D40" = "0.04232"
    "E40" = "1.50%"
    "D41" = "0.006731"
    "E41" = "1.06%"
    "E42" = "-12.21%"
    "E43" = "-8.35%"
    "D44" = "0.01286"
    "E44" = "3.60%"
    "D45" = "0.00005603"
    "E45" = "0.61%"
    "D46" = "1.973"
    "E46" = "26.43%"
    "E47" = "-29.41%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a plain-text number format first so Excel does not reinterpret
    # numeric- or percent-looking strings (e.g. "289.72", "1.01%") as real
    # numbers -- the source data is stored as text in this sheet.
    $cell.NumberFormat = "@"
    $cell.Value2 = $updates[$addr]
}
